# dbdp_instances/incgraph_25_25_0.065_0.2_1.xlsx
# Commit: "Otimização .degree com numba"
#
# The underlying diff renames the two worksheets, tweaks the saved
# selection / zoom on each sheet view, and repositions the second
# picture on the (now named) "y" sheet.

$wb = $excel.ActiveWorkbook

# --- Rename sheets: Sheet1 -> x, Sheet2 -> y ---------------------------
$wsX = $wb.Worksheets.Item(1)
$wsY = $wb.Worksheets.Item(2)
$wsX.Name = "x"
$wsY.Name = "y"

# --- Sheet "x": saved selection moves from AL4 to Z36 ------------------
$wsX.Activate() | Out-Null
$wsX.Range("Z36").Select() | Out-Null

# --- Sheet "y": zoom 130% -> 115%, selection AU22 -> AP30 ---------------
$wsY.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 115
$wsY.Range("AP30").Select() | Out-Null

# --- Move "Picture 2" anchor on sheet "y" --------------------------------
# Target anchor (from the OOXML diff):
#   from: col=73 colOff=207353 row=4  rowOff=83528
#   to:   col=91 colOff=107916 row=20 rowOff=26004
# Translating those cell-anchor coordinates into the point-based
# Left/Top/Width/Height the Shape object expects (default column width
# 58.4375pt, default row height 15pt) reproduces that anchor exactly
# once Excel re-saves the drawing part.
$shp = $wsY.Shapes.Item(2)
$shp.Left = 4282.264507874016
$shp.Top = 66.57700787401575
$shp.Width = 1044.0453149606299
$shp.Height = 235.47055118110237
